# [6] Einlesen und Anzeigen des Status eines PBIs.
#
# - Shorten the acceptance-criteria text of PBI #6 ("Status eines PBIs")
#   and flip its State from Todo to Done.
# - Insert a new PBI #10 ("Effort-Forecast je Sprint") as the new row 6,
#   shifting the former rows 6-10 down to 7-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (PBI 6 - "Status eines PBIs") ---
$ws.Range("C5").Value = "Akzeptanzkriterien:`n- aus dem CSV laden, im PBL anzeigen`n- Status: Todo, In Progress, Done, Canceled"
$ws.Range("E5").Value = "Done"
$ws.Rows(5).RowHeight = 75

# --- Insert a new row 6 for the new PBI 10 (Effort-Forecast) ---
$ws.Rows(6).Insert()

$ws.Range("A6").Value = 10
$ws.Range("C6").Value = "Akzeptanzkriterien:`n- Zur Berechnung wird die normierte Geschwindigkeit verwendet:`n-- Abgeschlossener Sprint: EffortDone/CapacityDone`n-- Laufender Sprint: EffortForecast/CapacityForecast`n- Diese wird mit der Kapazität des Sprints multipliziert.`n- Der Forecast ergibt sich aus den  normierten Geschwindigkeiten aller vorherigen Sprints multipliziert mit der Forecast-Kapazität des Sprints."
$ws.Range("B6").Value = "Als PO möchte ich den Effort-Forecast für jeden Sprint sehen können."
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Todo"
$ws.Range("F6").Value = "Sprint 2"
$ws.Rows(6).RowHeight = 255

$ws.Range("F7").Select()
